# Auto update Excel log
# Appends new sensor-log rows to the PIR, Humidity and Temperature sheets,
# mirroring the "Auto update Excel log" commit (2026-02-06 Bathroom readings).
#
# Every new cell is entered with a leading apostrophe so the value is stored
# as literal text (matching the source log's plain-text cells) instead of
# being auto-parsed into a date/time/percentage number by Excel's normal
# typed-input coercion. The quote-prefix formatting is then cleared by
# resetting the written range's Style back to "Normal" once all values are
# in place.

function Append-LogRows {
    param(
        [string]$SheetName,
        [int]$StartRow,
        [object[]]$Rows
    )

    $wb = $excel.ActiveWorkbook
    $ws = $wb.Worksheets.Item($SheetName)

    for ($i = 0; $i -lt $Rows.Count; $i++) {
        $r = $StartRow + $i
        $rowValues = $Rows[$i]
        for ($j = 0; $j -lt $rowValues.Count; $j++) {
            $col = $j + 1
            $ws.Cells.Item($r, $col).Value = "'" + $rowValues[$j]
        }
    }

    $endRow = $StartRow + $Rows.Count - 1
    $lastCol = $Rows[0].Count
    $colLetter = [char](64 + $lastCol)
    $rangeAddr = "A" + $StartRow + ":" + $colLetter + $endRow
    $ws.Range($rangeAddr).Style = "Normal"
}

# --- PIR sheet: append rows 259-270 ---------------------------------------
$pirRows = @(
    @("2026-02-06", "10:01:33", "10:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "10:01:35", "10:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "10:01:39", "10:00", "Bathroom", "Motion Detected", "Active"),
    @("2026-02-06", "10:01:46", "10:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "10:01:51", "10:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "10:01:56", "10:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "10:02:00", "10:00", "Bathroom", "Motion Detected", "Active"),
    @("2026-02-06", "10:02:07", "10:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "10:02:12", "10:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "10:02:17", "10:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "10:02:22", "10:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-06", "10:02:27", "10:00", "Bathroom", "No Motion", "Inactive")
)
Append-LogRows "PIR" 259 $pirRows

# --- Humidity sheet: append rows 162-172 ----------------------------------
$humidityRows = @(
    @("2026-02-06", "10:01:32", "10:00", "Bathroom", "68.3%", "Active"),
    @("2026-02-06", "10:01:34", "10:00", "Bathroom", "69.0%", "Active"),
    @("2026-02-06", "10:01:39", "10:00", "Bathroom", "68.0%", "Active"),
    @("2026-02-06", "10:01:44", "10:00", "Bathroom", "68.8%", "Active"),
    @("2026-02-06", "10:01:49", "10:00", "Bathroom", "68.8%", "Active"),
    @("2026-02-06", "10:01:54", "10:00", "Bathroom", "68.9%", "Active"),
    @("2026-02-06", "10:02:05", "10:00", "Bathroom", "68.9%", "Active"),
    @("2026-02-06", "10:02:09", "10:00", "Bathroom", "68.9%", "Active"),
    @("2026-02-06", "10:02:19", "10:00", "Bathroom", "68.1%", "Active"),
    @("2026-02-06", "10:02:24", "10:00", "Bathroom", "69.2%", "Active"),
    @("2026-02-06", "10:02:29", "10:00", "Bathroom", "68.5%", "Active")
)
Append-LogRows "Humidity" 162 $humidityRows

# --- Temperature sheet: append rows 162-172 -------------------------------
$temperatureRows = @(
    @("2026-02-06", "10:01:32", "10:00", "Bathroom", "27.9C", "Active"),
    @("2026-02-06", "10:01:35", "10:00", "Bathroom", "27.8C", "Active"),
    @("2026-02-06", "10:01:40", "10:00", "Bathroom", "28.0C", "Active"),
    @("2026-02-06", "10:01:45", "10:00", "Bathroom", "27.9C", "Active"),
    @("2026-02-06", "10:01:50", "10:00", "Bathroom", "27.9C", "Active"),
    @("2026-02-06", "10:01:55", "10:00", "Bathroom", "27.9C", "Active"),
    @("2026-02-06", "10:02:05", "10:00", "Bathroom", "27.9C", "Active"),
    @("2026-02-06", "10:02:10", "10:00", "Bathroom", "27.9C", "Active"),
    @("2026-02-06", "10:02:20", "10:00", "Bathroom", "27.8C", "Active"),
    @("2026-02-06", "10:02:25", "10:00", "Bathroom", "27.8C", "Active"),
    @("2026-02-06", "10:02:30", "10:00", "Bathroom", "27.9C", "Active")
)
Append-LogRows "Temperature" 162 $temperatureRows
